$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper-free, explicit cell writes.
# Rule learned from the runtime:
#   - to keep a cell a real NUMBER while still carrying the "Text" (@) number
#     format (matching style index 1 used throughout the data rows) the
#     Value must be assigned *before* NumberFormat is changed to "@".
#   - to force a numeric-looking string (e.g. "12", "12345", "2") to be
#     stored as TEXT, NumberFormat must be set to "@" *before* the Value is
#     assigned.
# ---------------------------------------------------------------------------

# --- fix existing row 2 values ---------------------------------------------
$ws.Range("G2").Value = "12"            # Month: December -> 12
$ws.Range("Q2").Value = "10"            # State: Indiana -> 10
$ws.Range("S2").Value = "21"            # Country: United States -> 21

$ws.Range("U2").NumberFormat = "@"
$ws.Range("U2").Value = "01234567890"   # Phone (text, leading zero)

$ws.Range("V2").NumberFormat = "@"
$ws.Range("V2").Value = "01787965432"   # Mobile (text, leading zero)

$ws.Range("R2").NumberFormat = "@"
$ws.Range("R2").Value = "12345"         # Zip: 1234 -> "12345" (text)

# --- add row 3 ---------------------------------------------------------------
$ws.Cells.Item(3,1).NumberFormat = "@"
$ws.Cells.Item(3,1).Value = "2"
$ws.Cells.Item(3,2).NumberFormat = "@"
$ws.Cells.Item(3,2).Value = "Ford"
$ws.Cells.Item(3,3).NumberFormat = "@"
$ws.Cells.Item(3,3).Value = "Prefect"
$ws.Cells.Item(3,4).NumberFormat = "@"
$ws.Cells.Item(3,4).Value = "aab@aab.com"
$ws.Cells.Item(3,5).NumberFormat = "@"
$ws.Cells.Item(3,5).Value = "revolution"
$ws.Cells.Item(3,6).Value = 31
$ws.Cells.Item(3,6).NumberFormat = "@"
$ws.Cells.Item(3,7).NumberFormat = "@"
$ws.Cells.Item(3,7).Value = "12"
$ws.Cells.Item(3,8).Value = 1983
$ws.Cells.Item(3,8).NumberFormat = "@"
$ws.Cells.Item(3,9).NumberFormat = "@"
$ws.Cells.Item(3,9).Value = "Y"
$ws.Cells.Item(3,10).NumberFormat = "@"
$ws.Cells.Item(3,10).Value = "Y"
$ws.Cells.Item(3,11).NumberFormat = "@"
$ws.Cells.Item(3,11).Value = "Karl"
$ws.Cells.Item(3,12).NumberFormat = "@"
$ws.Cells.Item(3,12).Value = "Marx"
$ws.Cells.Item(3,13).NumberFormat = "@"
$ws.Cells.Item(3,13).Value = "AAB"
$ws.Cells.Item(3,14).NumberFormat = "@"
$ws.Cells.Item(3,14).Value = "Trierer Stra" + [char]0xDF + "e"
$ws.Cells.Item(3,15).Value = 48
$ws.Cells.Item(3,15).NumberFormat = "@"
$ws.Cells.Item(3,16).NumberFormat = "@"
$ws.Cells.Item(3,16).Value = "New York"
$ws.Cells.Item(3,17).NumberFormat = "@"
$ws.Cells.Item(3,17).Value = "10"
$ws.Cells.Item(3,18).NumberFormat = "@"
$ws.Cells.Item(3,18).Value = "12345"
$ws.Cells.Item(3,19).NumberFormat = "@"
$ws.Cells.Item(3,19).Value = "21"
$ws.Cells.Item(3,20).NumberFormat = "@"
$ws.Cells.Item(3,20).Value = "Lorem ipsum dolor sit amet, consetetur sadipscing elitr, sed diam nonumy eirmod tempor invidunt ut labore et dolore magna aliquyam"
$ws.Cells.Item(3,21).NumberFormat = "@"
$ws.Cells.Item(3,21).Value = "01234567890"
$ws.Cells.Item(3,22).NumberFormat = "@"
$ws.Cells.Item(3,22).Value = "01787965432"
$ws.Cells.Item(3,23).NumberFormat = "@"
$ws.Cells.Item(3,23).Value = "thermiknator"

[void]$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:aab@aab.com")
$ws.Range("D3").Style = "Link"
$ws.Range("D3").NumberFormat = "@"

# --- add row 4 ---------------------------------------------------------------
$ws.Cells.Item(4,1).NumberFormat = "@"
$ws.Cells.Item(4,1).Value = "2"
$ws.Cells.Item(4,2).NumberFormat = "@"
$ws.Cells.Item(4,2).Value = "Ford"
$ws.Cells.Item(4,3).NumberFormat = "@"
$ws.Cells.Item(4,3).Value = "Prefect"
$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value = "aabtest1@aab.com"
$ws.Cells.Item(4,5).NumberFormat = "@"
$ws.Cells.Item(4,5).Value = "revolution"
$ws.Cells.Item(4,6).Value = 31
$ws.Cells.Item(4,6).NumberFormat = "@"
$ws.Cells.Item(4,7).NumberFormat = "@"
$ws.Cells.Item(4,7).Value = "12"
$ws.Cells.Item(4,8).Value = 1983
$ws.Cells.Item(4,8).NumberFormat = "@"
$ws.Cells.Item(4,9).NumberFormat = "@"
$ws.Cells.Item(4,9).Value = "Y"
$ws.Cells.Item(4,10).NumberFormat = "@"
$ws.Cells.Item(4,10).Value = "Y"
$ws.Cells.Item(4,11).NumberFormat = "@"
$ws.Cells.Item(4,11).Value = "Karl"
$ws.Cells.Item(4,12).NumberFormat = "@"
$ws.Cells.Item(4,12).Value = "Marx"
$ws.Cells.Item(4,13).NumberFormat = "@"
$ws.Cells.Item(4,13).Value = "AAB"
$ws.Cells.Item(4,14).NumberFormat = "@"
$ws.Cells.Item(4,14).Value = "Trierer Stra" + [char]0xDF + "e"
$ws.Cells.Item(4,15).Value = 48
$ws.Cells.Item(4,15).NumberFormat = "@"
$ws.Cells.Item(4,16).NumberFormat = "@"
$ws.Cells.Item(4,16).Value = "New York"
$ws.Cells.Item(4,17).NumberFormat = "@"
$ws.Cells.Item(4,17).Value = "10"
$ws.Cells.Item(4,18).NumberFormat = "@"
$ws.Cells.Item(4,18).Value = "12345"
$ws.Cells.Item(4,19).NumberFormat = "@"
$ws.Cells.Item(4,19).Value = "21"
$ws.Cells.Item(4,20).NumberFormat = "@"
$ws.Cells.Item(4,20).Value = "Lorem ipsum dolor sit amet, consetetur sadipscing elitr, sed diam nonumy eirmod tempor invidunt ut labore et dolore magna aliquyam"
$ws.Cells.Item(4,21).NumberFormat = "@"
$ws.Cells.Item(4,21).Value = "01234567890"
$ws.Cells.Item(4,22).NumberFormat = "@"
$ws.Cells.Item(4,22).Value = "01787965432"
$ws.Cells.Item(4,23).NumberFormat = "@"
$ws.Cells.Item(4,23).Value = "thermiknator"

[void]$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:aabtest1@aab.com")
$ws.Range("D4").Style = "Link"
$ws.Range("D4").NumberFormat = "@"

# --- view / selection state ---------------------------------------------------
[void]$ws.Range("E10").Select()
